# Actualizacion a 4 de Abril.
# Adds the new "Ecuador" column (I) and fills in the newly reported
# daily case numbers for the end of March / beginning of April.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# New column header
$ws.Range("I1").Value = "Ecuador"

# New "Ecuador" daily counts for rows that already had data
$ws.Range("I31").Value = 37
$ws.Range("I32").Value = 58
$ws.Range("I33").Value = 111
$ws.Range("I35").Value = 260
$ws.Range("I36").Value = 367
$ws.Range("I37").Value = 532
$ws.Range("I38").Value = 789
$ws.Range("I39").Value = 981
$ws.Range("I40").Value = 1082
$ws.Range("I41").Value = 1211
$ws.Range("I42").Value = 1403
$ws.Range("I44").Value = 1835
$ws.Range("I45").Value = 1924
$ws.Range("I46").Value = 1966
$ws.Range("I47").Value = 2302
$ws.Range("I49").Value = 3163
$ws.Range("I51").Value = 3465

# Newly reported data for rows 48-50 (previously only date/day were known)
$ws.Range("C48").Value = 3404
$ws.Range("D48").Value = 6880
$ws.Range("E48").Value = 104118
$ws.Range("F48").Value = 110574
$ws.Range("G48").Value = 9887
$ws.Range("H48").Value = 215003

$ws.Range("C49").Value = 3737
$ws.Range("D49").Value = 8044
$ws.Range("E49").Value = 112065
$ws.Range("F49").Value = 115242
$ws.Range("G49").Value = 9976
$ws.Range("H49").Value = 244877

$ws.Range("C50").Value = 4161
$ws.Range("D50").Value = 9194
$ws.Range("E50").Value = 119199
$ws.Range("F50").Value = 119827
$ws.Range("G50").Value = 10062
$ws.Range("H50").Value = 277161

# Update the view: keep gridlines visible, scroll the frozen pane down
# (rows unfroze/re-froze further down the sheet) and move the active
# selection to the new bottom-right cell (H51), matching the refreshed
# "latest data" viewport.
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("B22").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H51").Select()
